$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.272421479225159
$ws.Range("B1").Value = 1.783588528633118
$ws.Range("C1").Value = 3.213485479354858
$ws.Range("D1").Value = 3.809664487838745
$ws.Range("E1").Value = 1.257516264915466
